$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "sandbox" mode to "demo" mode.
$ws.Range("A69").Value = "demo"
$ws.Range("B69").Value = "DEMO"

# Remove the old "sandbox disclaimer" / "play anyway" / "play anyway gamepad" rows
# (rows 70-72), shifting everything below them up by three rows.
$ws.Range("A70:C72").EntireRow.Delete()

$ws.Range("A70").Select()
